# Edit script applying the "Natmi following Dr Hou advice" change:
# - Adds a new category "ECs" (EndoThelial Cells) as a sending/target cluster
# - Recomputes all LR-pair statistics rows for the 3x3 -> richer cluster grid
#   (M1, M2, Neutro, ECs) with updated cell counts, expression, and specificity values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "Ccl12"
$ws.Cells.Item(2, 3).Value = "Ccr5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 40.32940166666666
$ws.Cells.Item(2, 8).Value = 120.988205
$ws.Cells.Item(2, 9).Value = 0.4705770439863239
$ws.Cells.Item(2, 10).Value = 0.4705770439863239
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.03930266666666667
$ws.Cells.Item(2, 14).Value = 0.117908
$ws.Cells.Item(2, 15).Value = 0.0002899299635503584
$ws.Cells.Item(2, 16).Value = 0.0002899299635503584
$ws.Cells.Item(2, 17).Value = 1.585053030571111
$ws.Cells.Item(2, 18).Value = 14.26547727514
$ws.Cells.Item(2, 19).Value = 0.0001364343852105903
$ws.Cells.Item(2, 20).Value = 0.0001364343852105903

# Row 3
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "Ccl12"
$ws.Cells.Item(3, 3).Value = "Ccr5"
$ws.Cells.Item(3, 4).Value = "M1"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 40.32940166666666
$ws.Cells.Item(3, 8).Value = 120.988205
$ws.Cells.Item(3, 9).Value = 0.4705770439863239
$ws.Cells.Item(3, 10).Value = 0.4705770439863239
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 60.87605266666667
$ws.Cells.Item(3, 14).Value = 182.628158
$ws.Cells.Item(3, 15).Value = 0.4490736437918471
$ws.Cells.Item(3, 16).Value = 0.4490736437918471
$ws.Cells.Item(3, 17).Value = 2455.094779875154
$ws.Cells.Item(3, 18).Value = 22095.85301887639
$ws.Cells.Item(3, 19).Value = 0.2113237478277348
$ws.Cells.Item(3, 20).Value = 0.2113237478277348

# Row 4
$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "Ccl12"
$ws.Cells.Item(4, 3).Value = "Ccr5"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 40.32940166666666
$ws.Cells.Item(4, 8).Value = 120.988205
$ws.Cells.Item(4, 9).Value = 0.4705770439863239
$ws.Cells.Item(4, 10).Value = 0.4705770439863239
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 63.712864
$ws.Cells.Item(4, 14).Value = 191.138592
$ws.Cells.Item(4, 15).Value = 0.470000381752102
$ws.Cells.Item(4, 16).Value = 0.470000381752102
$ws.Cells.Item(4, 17).Value = 2569.501683589706
$ws.Cells.Item(4, 18).Value = 23125.51515230736
$ws.Cells.Item(4, 19).Value = 0.2211713903173479
$ws.Cells.Item(4, 20).Value = 0.2211713903173479

# Row 5
$ws.Cells.Item(5, 1).Value = "M1"
$ws.Cells.Item(5, 2).Value = "Ccl12"
$ws.Cells.Item(5, 3).Value = "Ccr5"
$ws.Cells.Item(5, 4).Value = "Neutro"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 40.32940166666666
$ws.Cells.Item(5, 8).Value = 120.988205
$ws.Cells.Item(5, 9).Value = 0.4705770439863239
$ws.Cells.Item(5, 10).Value = 0.4705770439863239
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 10.93095566666667
$ws.Cells.Item(5, 14).Value = 32.792867
$ws.Cells.Item(5, 15).Value = 0.08063604449250054
$ws.Cells.Item(5, 16).Value = 0.08063604449250053
$ws.Cells.Item(5, 17).Value = 440.8389016815261
$ws.Cells.Item(5, 18).Value = 3967.550115133735
$ws.Cells.Item(5, 19).Value = 0.0379454714560306
$ws.Cells.Item(5, 20).Value = 0.0379454714560306

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Ccl12"
$ws.Cells.Item(6, 3).Value = "Ccr5"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 35.65443533333333
$ws.Cells.Item(6, 8).Value = 106.963306
$ws.Cells.Item(6, 9).Value = 0.41602796200245
$ws.Cells.Item(6, 10).Value = 0.41602796200245
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.03930266666666667
$ws.Cells.Item(6, 14).Value = 0.117908
$ws.Cells.Item(6, 15).Value = 0.0002899299635503584
$ws.Cells.Item(6, 16).Value = 0.0002899299635503584
$ws.Cells.Item(6, 17).Value = 1.401314387094222
$ws.Cells.Item(6, 18).Value = 12.611829483848
$ws.Cells.Item(6, 19).Value = 0.0001206189718593002
$ws.Cells.Item(6, 20).Value = 0.0001206189718593002

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Ccl12"
$ws.Cells.Item(7, 3).Value = "Ccr5"
$ws.Cells.Item(7, 4).Value = "M1"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 35.65443533333333
$ws.Cells.Item(7, 8).Value = 106.963306
$ws.Cells.Item(7, 9).Value = 0.41602796200245
$ws.Cells.Item(7, 10).Value = 0.41602796200245
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 60.87605266666667
$ws.Cells.Item(7, 14).Value = 182.628158
$ws.Cells.Item(7, 15).Value = 0.4490736437918471
$ws.Cells.Item(7, 16).Value = 0.4490736437918471
$ws.Cells.Item(7, 17).Value = 2170.501283152261
$ws.Cells.Item(7, 18).Value = 19534.51154837035
$ws.Cells.Item(7, 19).Value = 0.1868271928157363
$ws.Cells.Item(7, 20).Value = 0.1868271928157363

# Row 8
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Ccl12"
$ws.Cells.Item(8, 3).Value = "Ccr5"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 35.65443533333333
$ws.Cells.Item(8, 8).Value = 106.963306
$ws.Cells.Item(8, 9).Value = 0.41602796200245
$ws.Cells.Item(8, 10).Value = 0.41602796200245
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 63.712864
$ws.Cells.Item(8, 14).Value = 191.138592
$ws.Cells.Item(8, 15).Value = 0.470000381752102
$ws.Cells.Item(8, 16).Value = 0.470000381752102
$ws.Cells.Item(8, 17).Value = 2271.646189389461
$ws.Cells.Item(8, 18).Value = 20444.81570450515
$ws.Cells.Item(8, 19).Value = 0.1955333009607005
$ws.Cells.Item(8, 20).Value = 0.1955333009607005

# Row 9
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Ccl12"
$ws.Cells.Item(9, 3).Value = "Ccr5"
$ws.Cells.Item(9, 4).Value = "Neutro"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 35.65443533333333
$ws.Cells.Item(9, 8).Value = 106.963306
$ws.Cells.Item(9, 9).Value = 0.41602796200245
$ws.Cells.Item(9, 10).Value = 0.41602796200245
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 10.93095566666667
$ws.Cells.Item(9, 14).Value = 32.792867
$ws.Cells.Item(9, 15).Value = 0.08063604449250054
$ws.Cells.Item(9, 16).Value = 0.08063604449250053
$ws.Cells.Item(9, 17).Value = 389.7370519487002
$ws.Cells.Item(9, 18).Value = 3507.633467538302
$ws.Cells.Item(9, 19).Value = 0.03354684925415388
$ws.Cells.Item(9, 20).Value = 0.03354684925415388

# Row 10
$ws.Cells.Item(10, 1).Value = "Neutro"
$ws.Cells.Item(10, 2).Value = "Ccl12"
$ws.Cells.Item(10, 3).Value = "Ccr5"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 9.718179666666666
$ws.Cells.Item(10, 8).Value = 29.154539
$ws.Cells.Item(10, 9).Value = 0.1133949940112261
$ws.Cells.Item(10, 10).Value = 0.1133949940112261
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.03930266666666667
$ws.Cells.Item(10, 14).Value = 0.117908
$ws.Cells.Item(10, 15).Value = 0.0002899299635503584
$ws.Cells.Item(10, 16).Value = 0.0002899299635503584
$ws.Cells.Item(10, 17).Value = 0.3819503760457778
$ws.Cells.Item(10, 18).Value = 3.437553384412
$ws.Cells.Item(10, 19).Value = [double]"3.28766064804679e-05"
$ws.Cells.Item(10, 20).Value = [double]"3.287660648046791e-05"

# Row 11
$ws.Cells.Item(11, 1).Value = "Neutro"
$ws.Cells.Item(11, 2).Value = "Ccl12"
$ws.Cells.Item(11, 3).Value = "Ccr5"
$ws.Cells.Item(11, 4).Value = "M1"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 9.718179666666666
$ws.Cells.Item(11, 8).Value = 29.154539
$ws.Cells.Item(11, 9).Value = 0.1133949940112261
$ws.Cells.Item(11, 10).Value = 0.1133949940112261
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 60.87605266666667
$ws.Cells.Item(11, 14).Value = 182.628158
$ws.Cells.Item(11, 15).Value = 0.4490736437918471
$ws.Cells.Item(11, 16).Value = 0.4490736437918471
$ws.Cells.Item(11, 17).Value = 591.6044172121292
$ws.Cells.Item(11, 18).Value = 5324.439754909163
$ws.Cells.Item(11, 19).Value = 0.05092270314837598
$ws.Cells.Item(11, 20).Value = 0.05092270314837598

# Row 12
$ws.Cells.Item(12, 1).Value = "Neutro"
$ws.Cells.Item(12, 2).Value = "Ccl12"
$ws.Cells.Item(12, 3).Value = "Ccr5"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 9.718179666666666
$ws.Cells.Item(12, 8).Value = 29.154539
$ws.Cells.Item(12, 9).Value = 0.1133949940112261
$ws.Cells.Item(12, 10).Value = 0.1133949940112261
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 63.712864
$ws.Cells.Item(12, 14).Value = 191.138592
$ws.Cells.Item(12, 15).Value = 0.470000381752102
$ws.Cells.Item(12, 16).Value = 0.470000381752102
$ws.Cells.Item(12, 17).Value = 619.1730594298987
$ws.Cells.Item(12, 18).Value = 5572.557534869088
$ws.Cells.Item(12, 19).Value = 0.05329569047405358
$ws.Cells.Item(12, 20).Value = 0.05329569047405359

# Row 13
$ws.Cells.Item(13, 1).Value = "Neutro"
$ws.Cells.Item(13, 2).Value = "Ccl12"
$ws.Cells.Item(13, 3).Value = "Ccr5"
$ws.Cells.Item(13, 4).Value = "Neutro"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 9.718179666666666
$ws.Cells.Item(13, 8).Value = 29.154539
$ws.Cells.Item(13, 9).Value = 0.1133949940112261
$ws.Cells.Item(13, 10).Value = 0.1133949940112261
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 10.93095566666667
$ws.Cells.Item(13, 14).Value = 32.792867
$ws.Cells.Item(13, 15).Value = 0.08063604449250054
$ws.Cells.Item(13, 16).Value = 0.08063604449250053
$ws.Cells.Item(13, 17).Value = 106.2289910970348
$ws.Cells.Item(13, 18).Value = 956.060919873313
$ws.Cells.Item(13, 19).Value = 0.00914372378231606
$ws.Cells.Item(13, 20).Value = 0.00914372378231606
